$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: "Fix: Can never reenter a server if you quit with TeamUI active" -> mark Completed = Yes
# (copy formatting from an existing "Completed = Yes" cell so it picks up the same
# highlighted style, then set the value)
$ws.Range("C16").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "Yes"

# Row 30: "Top down camera and camera rewrite" -> mark Completed = Yes, Completion Date = 2024-07-24 (serial 45497)
$ws.Range("C16:D16").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = "Yes"
$ws.Range("D30").Value = 45497

$excel.CutCopyMode = 0

# Update viewport / selection to match final state
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A39").Select()
